$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update path to video file (D2) with the new file name
$ws.Range("D2").Value = "C:\Phase1\Videos\Men Will Be Men - edited.mp4"

# Update wait time (sec) after video ends (E2) - now 77 seconds (0:01:17)
$ws.Range("E2").Value = 0.00089120370370370362

# Update the selected cell in the sheet view
$null = $ws.Range("E3").Select()
